$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 441.66666
$ws.Range("J4").Value = 733.3333
$ws.Range("L4").Value = 733.3333
$ws.Range("N4").Value = -961.3333
$ws.Range("H62").Value = 3339.25
$ws.Range("I62").Value = 2663.5
$ws.Range("K62").Value = 2663.5
$ws.Range("M62").Value = -2039.5
$ws.Range("H65").Value = 3339.25
$ws.Range("I65").Value = 2663.5
$ws.Range("K65").Value = 13317.5
$ws.Range("M65").Value = -10197.5
$ws.Range("H86").Value = 20105
$ws.Range("I86").Value = 10271.857
$ws.Range("K86").Value = 10271.857
$ws.Range("M86").Value = -9148.857
$ws.Range("H89").Value = 20105
$ws.Range("I89").Value = 10271.857
$ws.Range("K89").Value = 51359.285
$ws.Range("M89").Value = -45743.285
$ws.Range("H137").Value = 20473.396
$ws.Range("I137").Value = 1387.3334
$ws.Range("J137").Value = 73641.71000000001
$ws.Range("K137").Value = 4162.0002
$ws.Range("L137").Value = 220925.13
$ws.Range("M137").Value = -1612.0002
$ws.Range("N137").Value = -226025.13
$ws.Range("H138").Value = 14708686
$ws.Range("I138").Value = 200005070
$ws.Range("J138").Value = 2623.4443
$ws.Range("K138").Value = 600015210
$ws.Range("L138").Value = 7870.3329
$ws.Range("M138").Value = -600010070
$ws.Range("N138").Value = -18150.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("M4").Value = -134
$ws.Range("H5").Value = 95
$ws.Range("I5").Value = 95
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 95
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 17
$ws.Range("N5").ClearContents()
$ws.Range("H21").Value = 2457.5
$ws.Range("I21").Value = 2457.5
$ws.Range("K21").Value = 2457.5
$ws.Range("M21").Value = -2083.5
$ws.Range("H30").Value = 1000
$ws.Range("J30").Value = 1000
$ws.Range("L30").Value = 1000
$ws.Range("N30").Value = -1300
$ws.Range("H45").Value = 3459
$ws.Range("I45").Value = 3214.0625
$ws.Range("K45").Value = 3214.0625
$ws.Range("M45").Value = -2837.0625
$ws.Range("H61").Value = 1611.2424
$ws.Range("I61").Value = 1264.3928
$ws.Range("K61").Value = 1264.3928
$ws.Range("M61").Value = -1052.3928
$ws.Range("H102").Value = 1750
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 122
$ws.Range("N102").Value = -5244
$ws.Range("H110").Value = 473.47058
$ws.Range("I110").Value = 489.5
$ws.Range("J110").Value = 435
$ws.Range("K110").Value = 489.5
$ws.Range("L110").Value = 435
$ws.Range("M110").Value = 1555.5
$ws.Range("N110").Value = -4525
$ws.Range("H122").Value = 1833.5294
$ws.Range("I122").Value = 1505
$ws.Range("K122").Value = 4515
$ws.Range("M122").Value = -2065
$ws.Range("H132").Value = 40552.54
$ws.Range("I132").Value = 2133.8235
$ws.Range("J132").Value = 301799.8
$ws.Range("K132").Value = 6401.470499999999
$ws.Range("L132").Value = 905399.3999999999
$ws.Range("M132").Value = -3871.470499999999
$ws.Range("N132").Value = -910459.3999999999
$ws.Range("H136").Value = 1611.2424
$ws.Range("I136").Value = 1264.3928
$ws.Range("K136").Value = 3793.1784
$ws.Range("M136").Value = -1243.1784

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 95
$ws.Range("I4").Value = 95
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 95
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 20
$ws.Range("N4").ClearContents()
$ws.Range("H94").Value = 750
$ws.Range("I94").Value = 750
$ws.Range("K94").Value = 750
$ws.Range("M94").Value = -299
$ws.Range("H99").Value = 1336.1177
$ws.Range("J99").Value = 1250.6
$ws.Range("L99").Value = 1250.6
$ws.Range("N99").Value = -4246.6
$ws.Range("H105").Value = 1658.5652
$ws.Range("I105").Value = 1176.4667
$ws.Range("K105").Value = 1176.4667
$ws.Range("M105").Value = 570.5333000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 312.52
$ws.Range("I22").Value = 310
$ws.Range("J22").Value = 315.72726
$ws.Range("K22").Value = 310
$ws.Range("L22").Value = 315.72726
$ws.Range("M22").Value = 40
$ws.Range("N22").Value = -1015.72726
$ws.Range("H122").Value = 1626.5186
$ws.Range("I122").Value = 1834.6923
$ws.Range("J122").Value = 1433.2142
$ws.Range("K122").Value = 5504.0769
$ws.Range("L122").Value = 4299.642599999999
$ws.Range("M122").Value = -3054.0769
$ws.Range("N122").Value = -9199.642599999999
$ws.Range("H132").Value = 8956.282999999999
$ws.Range("I132").Value = 9469.779
$ws.Range("K132").Value = 28409.337
$ws.Range("M132").Value = -25879.337
$ws.Range("H134").Value = 660.69385
$ws.Range("I134").Value = 548.3077
$ws.Range("J134").Value = 1099
$ws.Range("K134").Value = 1644.9231
$ws.Range("L134").Value = 3297
$ws.Range("M134").Value = 890.0769
$ws.Range("N134").Value = -8367
$ws.Range("H138").Value = 49939.09
$ws.Range("J138").Value = 49939.09
$ws.Range("L138").Value = 49939.09
$ws.Range("N138").Value = -60219.09

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 5904.5884
$ws.Range("I2").Value = 7707.077
$ws.Range("J2").Value = 46.5
$ws.Range("K2").Value = 46242.462
$ws.Range("L2").Value = 279
$ws.Range("M2").Value = -46129.462
$ws.Range("N2").Value = -505
$ws.Range("H38").Value = 62500080
$ws.Range("J38").Value = 250000030
$ws.Range("L38").Value = 750000090
$ws.Range("N38").Value = -750000784
$ws.Range("H131").Value = 744.0599999999999
$ws.Range("J131").Value = 793.1429000000001
$ws.Range("L131").Value = 2379.4287
$ws.Range("N131").Value = -12459.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 44502.805
$ws.Range("I132").Value = 38671.535
$ws.Range("J132").Value = 64912.25
$ws.Range("K132").Value = 116014.605
$ws.Range("L132").Value = 194736.75
$ws.Range("M132").Value = -113484.605
$ws.Range("N132").Value = -199796.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 823.8889
$ws.Range("I46").Value = 866.1818
$ws.Range("K46").Value = 866.1818
$ws.Range("M46").Value = -678.1818
$ws.Range("H93").Value = 1126.9412
$ws.Range("I93").Value = 1054.1428
$ws.Range("K93").Value = 1054.1428
$ws.Range("M93").Value = 193.8571999999999
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H132").Value = 2884.1667
$ws.Range("I132").Value = 2527
$ws.Range("K132").Value = 7581
$ws.Range("M132").Value = -5051

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 7490
$ws.Range("J31").Value = 7490
$ws.Range("L31").Value = 7490
$ws.Range("N31").Value = -8186
$ws.Range("H132").Value = 2108.4167
$ws.Range("I132").Value = 1533.6666
$ws.Range("J132").Value = 3832.6667
$ws.Range("K132").Value = 4600.9998
$ws.Range("L132").Value = 11498.0001
$ws.Range("M132").Value = -2070.9998
$ws.Range("N132").Value = -16558.0001
